# "Addressing mode" switcher -> select:
# adds a new "Addressing mode(Optional)" column (O) to the node-info
# template, with "Static"/"Dynamic" choices for the two example rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O: header first, then the data rows. Row 3 ("Static") is
# written before row 2 ("Dynamic") so the new entries land in the shared
# string table in the same order as the authored workbook (Addressing
# mode(Optional), Static, Dynamic).
$ws.Range("O1").Value = "Addressing mode(Optional)"
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Match the widened column O from the source edit (closest value this
# engine's pixel-quantised ColumnWidth can reach to the authored 27.125
# character-width: it rounds to the nearest 1/7th of a character).
$ws.Columns.Item(15).ColumnWidth = 26.428571428571427

# The existing "Speed limit M/s(Optional)" header (N1) also picked up a
# distinct font on its zero-width-space characters in the source edit;
# reproduce that run split here.
$n1 = $ws.Cells.Item(1, 14)
$zwsp = $n1.Characters(7, 2)
$zwsp.Font.Name = "MS Gothic"
$zwsp.Font.Size = 12

# Author's selection at save time.
$null = $ws.Range("G17").Select()
